$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 92756.17999999999
$ws.Range("I4").Value = 143910.58
$ws.Range("J4").Value = 3236
$ws.Range("K4").Value = 143910.58
$ws.Range("L4").Value = 3236
$ws.Range("M4").Value = -143796.58
$ws.Range("N4").Value = -3464
$ws.Range("H5").Value = 113.625
$ws.Range("I5").Value = 96.833336
$ws.Range("K5").Value = 96.833336
$ws.Range("M5").Value = 18.166664
$ws.Range("H19").Value = 1060.75
$ws.Range("J19").Value = 1214.5
$ws.Range("L19").Value = 1214.5
$ws.Range("N19").Value = -1564.5
$ws.Range("H29").Value = 4214.2856
$ws.Range("J29").Value = 4833.3335
$ws.Range("L29").Value = 14500.0005
$ws.Range("N29").Value = -15062.0005
$ws.Range("H32").Value = 27779096
$ws.Range("I32").Value = 946.1429000000001
$ws.Range("J32").Value = 45456100
$ws.Range("K32").Value = 946.1429000000001
$ws.Range("L32").Value = 45456100
$ws.Range("M32").Value = -620.1429000000001
$ws.Range("N32").Value = -45456752
$ws.Range("H41").Value = 709.2
$ws.Range("J41").Value = 1062.125
$ws.Range("L41").Value = 1062.125
$ws.Range("N41").Value = -1942.125
$ws.Range("H43").Value = 48149596
$ws.Range("I43").Value = 72223210
$ws.Range("K43").Value = 72223210
$ws.Range("M43").Value = -72223141
$ws.Range("H53").Value = 200001360
$ws.Range("I53").Value = 1650
$ws.Range("K53").Value = 1650
$ws.Range("M53").Value = -1013
$ws.Range("H74").Value = 5677.875
$ws.Range("I74").Value = 5789.7334
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 5789.7334
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -4853.7334
$ws.Range("N74").Value = -5872
$ws.Range("H76").Value = 3489.25
$ws.Range("I76").Value = 3489.25
$ws.Range("K76").Value = 3489.25
$ws.Range("M76").Value = -3174.25
$ws.Range("H77").Value = 5677.875
$ws.Range("I77").Value = 5789.7334
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 28948.667
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -24268.667
$ws.Range("N77").Value = -29360
$ws.Range("H79").Value = 3489.25
$ws.Range("I79").Value = 3489.25
$ws.Range("K79").Value = 3489.25
$ws.Range("M79").Value = -2397.25
$ws.Range("H132").Value = 1796.0344
$ws.Range("I132").Value = 1767.4
$ws.Range("K132").Value = 5302.200000000001
$ws.Range("M132").Value = -2772.200000000001
$ws.Range("H137").Value = 4009374.5
$ws.Range("I137").Value = 7020.4443
$ws.Range("J137").Value = 14301142
$ws.Range("K137").Value = 21061.3329
$ws.Range("L137").Value = 42903426
$ws.Range("M137").Value = -18511.3329
$ws.Range("N137").Value = -42908526
$ws.Range("H138").Value = 5712.1763
$ws.Range("I138").Value = 6489.6924
$ws.Range("K138").Value = 19469.0772
$ws.Range("M138").Value = -14329.0772
$ws.Range("H140").Value = 92217.8
$ws.Range("J140").Value = 94578.164
$ws.Range("L140").Value = 94578.164
$ws.Range("N140").Value = -104938.164

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 519
$ws.Range("I4").Value = 873.6667
$ws.Range("K4").Value = 873.6667
$ws.Range("M4").Value = -757.6667
$ws.Range("H12").Value = 1999.5
$ws.Range("I12").Value = 1999.5
$ws.Range("K12").Value = 1999.5
$ws.Range("M12").Value = -1826.5
$ws.Range("H32").Value = 168900.83
$ws.Range("I32").Value = 193346.77
$ws.Range("J32").Value = 10002.25
$ws.Range("K32").Value = 193346.77
$ws.Range("L32").Value = 10002.25
$ws.Range("M32").Value = -193059.77
$ws.Range("N32").Value = -10576.25
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H63").Value = 13574.777
$ws.Range("J63").Value = 17158.814
$ws.Range("L63").Value = 17158.814
$ws.Range("N63").Value = -18530.814
$ws.Range("H66").Value = 13574.777
$ws.Range("J66").Value = 17158.814
$ws.Range("L66").Value = 85794.06999999999
$ws.Range("N66").Value = -92658.06999999999
$ws.Range("H74").Value = 449003.78
$ws.Range("I74").Value = 1902.5264
$ws.Range("J74").Value = 1662564.2
$ws.Range("K74").Value = 1902.5264
$ws.Range("L74").Value = 1662564.2
$ws.Range("M74").Value = -1028.5264
$ws.Range("N74").Value = -1664312.2
$ws.Range("H77").Value = 449003.78
$ws.Range("I77").Value = 1902.5264
$ws.Range("J77").Value = 1662564.2
$ws.Range("K77").Value = 9512.632
$ws.Range("L77").Value = 8312821
$ws.Range("M77").Value = -5144.632
$ws.Range("N77").Value = -8321557
$ws.Range("H102").Value = 2245
$ws.Range("I102").Value = 1991
$ws.Range("K102").Value = 1991
$ws.Range("M102").Value = -369
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("H132").Value = 3088.0527
$ws.Range("I132").Value = 2450.5
$ws.Range("J132").Value = 3551.7273
$ws.Range("K132").Value = 7351.5
$ws.Range("L132").Value = 10655.1819
$ws.Range("M132").Value = -4821.5
$ws.Range("N132").Value = -15715.1819

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1406.8276
$ws.Range("I94").Value = 1522.9565
$ws.Range("K94").Value = 1522.9565
$ws.Range("M94").Value = -1071.9565
$ws.Range("H105").Value = 11448.381
$ws.Range("I105").Value = 14363.777
$ws.Range("J105").Value = 9261.833000000001
$ws.Range("K105").Value = 14363.777
$ws.Range("L105").Value = 9261.833000000001
$ws.Range("M105").Value = -12616.777
$ws.Range("N105").Value = -12755.833

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 174.7
$ws.Range("I7").Value = 188.55556
$ws.Range("K7").Value = 188.55556
$ws.Range("M7").Value = -75.55556000000001
$ws.Range("H10").Value = 3196.4
$ws.Range("I10").Value = 3171
$ws.Range("J10").Value = 3234.5
$ws.Range("K10").Value = 3171
$ws.Range("L10").Value = 3234.5
$ws.Range("M10").Value = -3032
$ws.Range("N10").Value = -3512.5
$ws.Range("H22").Value = 541.5122
$ws.Range("I22").Value = 541.5122
$ws.Range("K22").Value = 541.5122
$ws.Range("M22").Value = -191.5122
$ws.Range("H31").Value = 4463.8276
$ws.Range("I31").Value = 3770.077
$ws.Range("J31").Value = 5027.5
$ws.Range("K31").Value = 3770.077
$ws.Range("L31").Value = 5027.5
$ws.Range("M31").Value = -3475.077
$ws.Range("N31").Value = -5617.5
$ws.Range("H34").Value = 4463.8276
$ws.Range("I34").Value = 3770.077
$ws.Range("J34").Value = 5027.5
$ws.Range("K34").Value = 3770.077
$ws.Range("L34").Value = 5027.5
$ws.Range("M34").Value = -3568.077
$ws.Range("N34").Value = -5431.5
$ws.Range("H44").Value = 5484.2856
$ws.Range("J44").Value = 5898.3335
$ws.Range("L44").Value = 5898.3335
$ws.Range("N44").Value = -6782.3335
$ws.Range("H86").Value = 12968.526
$ws.Range("I86").Value = 5436.5454
$ws.Range("K86").Value = 5436.5454
$ws.Range("M86").Value = -4313.5454
$ws.Range("H89").Value = 12968.526
$ws.Range("I89").Value = 5436.5454
$ws.Range("K89").Value = 27182.727
$ws.Range("M89").Value = -21566.727
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H132").Value = 1758.1666
$ws.Range("J132").Value = 2226.2856
$ws.Range("L132").Value = 6678.8568
$ws.Range("N132").Value = -11738.8568
$ws.Range("H141").Value = 401707.66
$ws.Range("I141").Value = 949999
$ws.Range("J141").Value = 333171.25
$ws.Range("K141").Value = 949999
$ws.Range("L141").Value = 333171.25
$ws.Range("M141").Value = -944819
$ws.Range("N141").Value = -343531.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 614.25
$ws.Range("I5").Value = 636.86664
$ws.Range("J5").Value = 275
$ws.Range("K5").Value = 1910.59992
$ws.Range("L5").Value = 825
$ws.Range("M5").Value = -1798.59992
$ws.Range("N5").Value = -1049
$ws.Range("H11").Value = 529.2727
$ws.Range("J11").Value = 477.4
$ws.Range("L11").Value = 1432.2
$ws.Range("N11").Value = -1712.2
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("H25").Value = 3249.75
$ws.Range("I25").Value = 999.6667
$ws.Range("K25").Value = 2999.0001
$ws.Range("M25").Value = -2830.0001
$ws.Range("H30").Value = 3249.75
$ws.Range("I30").Value = 999.6667
$ws.Range("K30").Value = 2999.0001
$ws.Range("M30").Value = -2897.0001
$ws.Range("H56").Value = 7699548.5
$ws.Range("I56").Value = 7699548.5
$ws.Range("K56").Value = 7699548.5
$ws.Range("M56").Value = -7699018.5
$ws.Range("H122").Value = 5129340
$ws.Range("I122").Value = 6667491
$ws.Range("K122").Value = 60007419
$ws.Range("M122").Value = -60004969
$ws.Range("H135").Value = 614.25
$ws.Range("I135").Value = 636.86664
$ws.Range("J135").Value = 275
$ws.Range("K135").Value = 5731.79976
$ws.Range("L135").Value = 2475
$ws.Range("M135").Value = -3196.79976
$ws.Range("N135").Value = -7545
$ws.Range("H141").Value = 2680.4443
$ws.Range("I141").Value = 2640.5
$ws.Range("K141").Value = 7921.5
$ws.Range("M141").Value = -2741.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 95.666664
$ws.Range("I2").Value = 94.22221999999999
$ws.Range("K2").Value = 94.22221999999999
$ws.Range("M2").Value = 18.77778000000001
$ws.Range("H3").Value = 801358.8
$ws.Range("J3").Value = 2002000
$ws.Range("L3").Value = 2002000
$ws.Range("N3").Value = -2002232
$ws.Range("H10").Value = 3489844.5
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H11").Value = 2417000
$ws.Range("I11").Value = 28333.334
$ws.Range("K11").Value = 28333.334
$ws.Range("M11").Value = -28194.334
$ws.Range("H18").Value = 30000
$ws.Range("I18").Value = 30000
$ws.Range("K18").Value = 30000
$ws.Range("M18").Value = -29707
$ws.Range("H24").Value = 23000
$ws.Range("J24").Value = 23000
$ws.Range("L24").Value = 23000
$ws.Range("N24").Value = -23346
$ws.Range("H35").Value = 12000
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 12000
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 12000
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -12596
$ws.Range("H122").Value = 2384.6296
$ws.Range("I122").Value = 2016
$ws.Range("J122").Value = 3437.8572
$ws.Range("K122").Value = 6048
$ws.Range("L122").Value = 10313.5716
$ws.Range("M122").Value = -3598
$ws.Range("N122").Value = -15213.5716

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 12077
$ws.Range("I7").Value = 7201.6
$ws.Range("J7").Value = 15124.125
$ws.Range("K7").Value = 7201.6
$ws.Range("L7").Value = 15124.125
$ws.Range("M7").Value = -7089.6
$ws.Range("N7").Value = -15348.125
$ws.Range("H18").Value = 38777.777
$ws.Range("I18").Value = 34500
$ws.Range("K18").Value = 34500
$ws.Range("M18").Value = -34328
$ws.Range("H22").Value = 3353
$ws.Range("I22").Value = 806.8570999999999
$ws.Range("J22").Value = 4466.9375
$ws.Range("K22").Value = 806.8570999999999
$ws.Range("L22").Value = 4466.9375
$ws.Range("M22").Value = -511.8570999999999
$ws.Range("N22").Value = -5056.9375
$ws.Range("H23").Value = 1672666.6
$ws.Range("I23").Value = 5000000
$ws.Range("K23").Value = 5000000
$ws.Range("M23").Value = -4999770
$ws.Range("H27").Value = 3353
$ws.Range("I27").Value = 806.8570999999999
$ws.Range("J27").Value = 4466.9375
$ws.Range("K27").Value = 806.8570999999999
$ws.Range("L27").Value = 4466.9375
$ws.Range("M27").Value = -699.8570999999999
$ws.Range("N27").Value = -4680.9375
$ws.Range("H33").Value = 10000
$ws.Range("I33").Value = 10000
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 10000
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -9710
$ws.Range("N33").ClearContents()
$ws.Range("H43").Value = 2011592.9
$ws.Range("I43").Value = 142400
$ws.Range("K43").Value = 142400
$ws.Range("M43").Value = -142207
$ws.Range("H46").Value = 7388.0835
$ws.Range("I46").Value = 35999.332
$ws.Range("J46").Value = 3300.762
$ws.Range("K46").Value = 35999.332
$ws.Range("L46").Value = 3300.762
$ws.Range("M46").Value = -35811.332
$ws.Range("N46").Value = -3676.762
$ws.Range("H82").Value = 1249.5
$ws.Range("I82").Value = 1249.5
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 1249.5
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -888.5
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 1249.5
$ws.Range("I85").Value = 1249.5
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 1249.5
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -1.5
$ws.Range("N85").ClearContents()
$ws.Range("H106").Value = 17456.334
$ws.Range("J106").Value = 17456.334
$ws.Range("L106").Value = 17456.334
$ws.Range("N106").Value = -19980.334
$ws.Range("H126").Value = 12077
$ws.Range("I126").Value = 7201.6
$ws.Range("J126").Value = 15124.125
$ws.Range("K126").Value = 21604.8
$ws.Range("L126").Value = 45372.375
$ws.Range("M126").Value = -19134.8
$ws.Range("N126").Value = -50312.375
$ws.Range("H132").Value = 5993.125
$ws.Range("I132").Value = 5325
$ws.Range("J132").Value = 6394
$ws.Range("K132").Value = 15975
$ws.Range("L132").Value = 19182
$ws.Range("M132").Value = -13445
$ws.Range("N132").Value = -24242

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 253749.75
$ws.Range("I81").Value = 5000
$ws.Range("K81").Value = 10000
$ws.Range("M81").Value = -8939
$ws.Range("H84").Value = 253749.75
$ws.Range("I84").Value = 5000
$ws.Range("K84").Value = 50000
$ws.Range("M84").Value = -44696
$ws.Range("H104").Value = 30576.4
$ws.Range("J104").Value = 30576.4
$ws.Range("L104").Value = 30576.4
$ws.Range("N104").Value = -37564.4
$ws.Range("H107").Value = 1145188.9
$ws.Range("I107").Value = 2214.6667
$ws.Range("J107").Value = 2200242
$ws.Range("K107").Value = 6644.000100000001
$ws.Range("L107").Value = 6600726
$ws.Range("M107").Value = -4724.000100000001
$ws.Range("N107").Value = -6604566
